# Generate Report for Handback
# The file 95e730fc-4864-4976-9d13-5bca16481286.md (zh-cn and de-de) has
# moved from "Ready for handoff" to "Handed back: in sync with en-US".
# Update the Status on the Overview sheet plus each language sheet, and
# stamp the Latest Handback DateTime for each language.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-02-24 12:03:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-02-24 12:03:39"
